$d = $word.ActiveDocument

# Locate the "| Internal circulation only" heading paragraph. The new
# "{{date}} " run must be inserted as the very first run of that
# paragraph (immediately after the paragraph properties, before the
# existing "|" run).
$target = $d.Content
$found = $target.Find.Execute("| Internal circulation only", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    # Fallback: try matching on just the distinctive trailing phrase in
    # case leading characters/spacing differ slightly.
    $target = $d.Content
    $found = $target.Find.Execute("Internal circulation only", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

if ($found) {
    $insertionPoint = $d.Range($target.Start, $target.Start)

    $xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Gill Sans Nova" w:hAnsi="Gill Sans Nova"/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="darkBlue"/></w:rPr><w:t xml:space="preserve">{{date}} </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part>
</pkg:package>
'@

    $insertionPoint.InsertXML($xml)
    Write-Host "Inserted {{date}} placeholder run before the 'Internal circulation only' heading."
} else {
    Write-Host "ERROR: could not locate target paragraph for {{date}} insertion."
}
